$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full refreshed cryptocurrency data table (rows 2-51, columns B-G)
# Column A (Rank) is unchanged (1-50)
$data = @(
    @('BTC', 'Bitcoin', 63715, 1255762952937, 17370474932, -0.29686),
    @('ETH', 'Ethereum', 3130.78, 376603418229, 8763130860, -0.28135),
    @('USDT', 'Tether', 0.999902, 110873401630, 25225379420, -0.02078),
    @('BNB', 'BNB', 585.74, 90173064375, 635538159, -0.77922),
    @('SOL', 'Solana', 146.05, 65529412496, 3131231336, -0.789),
    @('USDC', 'USDC', 0.99929, 33571337848, 5520914577, -0.08293),
    @('STETH', 'Lido Staked Ether', 3131.44, 29347777369, 48847870, -0.19961),
    @('XRP', 'XRP', 0.528321, 29221706776, 553316388, -1.29514),
    @('DOGE', 'Dogecoin', 0.160376, 23139460998, 2522810945, 0.3014),
    @('TON', 'Toncoin', 5.75, 20046557385, 177504637, 0.17876),
    @('ADA', 'Cardano', 0.456519, 16163549142, 250731059, -2.83631),
    @('SHIB', 'Shiba Inu', 0.00002466, 14570427080, 562756528, -3.40666),
    @('AVAX', 'Avalanche', 36.89, 14031422692, 495062915, 2.57879),
    @('TRX', 'TRON', 0.121154, 10616002690, 250269645, -1.62527),
    @('WBTC', 'Wrapped Bitcoin', 63614, 9885751780, 153718717, -0.41951),
    @('DOT', 'Polkadot', 7.05, 9614535142, 148223024, -1.57791),
    @('BCH', 'Bitcoin Cash', 462.83, 9136636752, 230303298, -2.23545),
    @('LINK', 'Chainlink', 14.26, 8383312870, 235240482, 0.2694),
    @('NEAR', 'NEAR Protocol', 7.1, 7615438426, 286576740, 3.86795),
    @('MATIC', 'Polygon', 0.732311, 6810152586, 187790863, 0.03065),
    @('LTC', 'Litecoin', 80.89, 6038089372, 244743312, -1.7276),
    @('ICP', 'Internet Computer', 12.89, 5985213594, 108591653, -3.53464),
    @('UNI', 'Uniswap', 7.41, 5604865212, 109649934, -1.82244),
    @('FET', 'Fetch.ai', 2.2, 5570010786, 233439366, 0.34556),
    @('DAI', 'Dai', 0.999424, 5439942790, 500299494, -0.06643),
    @('LEO', 'LEO Token', 5.76, 5321884361, 1036789, -2.04791),
    @('ETC', 'Ethereum Classic', 26.83, 3954726449, 138473074, -0.95644),
    @('HBAR', 'Hedera', 0.108682, 3895222731, 99877542, -0.23268),
    @('APT', 'Aptos', 9.04, 3870674109, 100797837, -1.52464),
    @('FDUSD', 'First Digital USD', 0.999702, 3851112521, 3420858087, 0.00936),
    @('RNDR', 'Render', 9.31, 3612223328, 222383049, 7.28763),
    @('PEPE', 'Pepe', 0.00000849, 3568422236, 794437951, -3.24556),
    @('CRO', 'Cronos', 0.131978, 3531686555, 8920661, -2.00962),
    @('ATOM', 'Cosmos Hub', 8.85, 3463237811, 164304035, 0.81311),
    @('MNT', 'Mantle', 1.042, 3404435412, 45617338, -1.05121),
    @('STX', 'Stacks', 2.31, 3381979948, 65097875, -4.59003),
    @('WIF', 'dogwifhat', 3.34, 3327025990, 577580828, -1.91548),
    @('FIL', 'Filecoin', 6, 3293020844, 121636321, -2.27053),
    @('IMX', 'Immutable', 2.21, 3225875891, 40849830, -0.39315),
    @('XLM', 'Stellar', 0.110094, 3190169472, 42863956, -1.41801),
    @('WEETH', 'Wrapped eETH', 3245.74, 3121062474, 37272846, -0.15604),
    @('EZETH', 'Renzo Restaked ETH', 3083.09, 3078346946, 26804255, -0.02081),
    @('OKB', 'OKB', 51.18, 3071344610, 5749835, 0.57407),
    @('OP', 'Optimism', 2.89, 3023067081, 262692274, -1.79478),
    @('XT', 'XT.com', 2.98, 2980539272, 659540, -5.48551),
    @('TAO', 'Bittensor', 435.37, 2932471137, 28005813, -2.2754),
    @('ARB', 'Arbitrum', 1.05, 2787204227, 216046746, -2.09126),
    @('MKR', 'Maker', 2909.65, 2698000940, 71951152, -0.54273),
    @('VET', 'VeChain', 0.03693267, 2692056731, 46001947, -0.77693),
    @('GRT', 'The Graph', 0.278586, 2651767808, 85567902, -1.9587)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 2).Value = $rowVals[0]
    $ws.Cells.Item($r, 3).Value = $rowVals[1]
    $ws.Cells.Item($r, 4).Value = $rowVals[2]
    $ws.Cells.Item($r, 5).Value = $rowVals[3]
    $ws.Cells.Item($r, 6).Value = $rowVals[4]
    $ws.Cells.Item($r, 7).Value = $rowVals[5]
}
